$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ID_0001")
$ws.Activate()

# Delete the second row (which contained the automation.devmrkolv@gmail.com /
# $chlUe13elKiNd credentials + its hyperlink), shifting the rows below it up.
$ws.Rows.Item(2).Delete()

# The former row 3 (now row 2) has its B cell set to mirror A (both "mmm").
$ws.Range("B2").Value2 = $ws.Range("A2").Value2

# Remove the now-stale hyperlink definition that pointed at the deleted row.
$ws.Hyperlinks.Delete()

# Leave the selection on B2, as last left by the editor.
$ws.Range("B2").Select()
